# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2024-10-17 (serial 45582) to 2024-10-18 (serial 45583).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45582) {
        $cell.Value2 = 45583
    }
}
